$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The task list is sorted by Status (column D) descending; newly-finished
# items move up next to the other "Done" rows. Rather than fight the grid's
# own re-sort, place every task at its final resting row directly.

# Row 2: clarify that the pop-up only happens on Vodafone
$ws.Range("A2").Value = "Lidar com pop-up menu após chamada (só vodafone?)"
$ws.Range("B2").Value = 1
$ws.Range("D2").Value = "TBD"
$ws.Range("F2").Value = "new functionality"

# Row 3: unchanged
$ws.Range("A3").Value = "Always visible"
$ws.Range("B3").Value = 1
$ws.Range("D3").Value = "TBD"
$ws.Range("F3").Value = "new functionality"
$ws.Range("G3").Value = "what if other activity from another app shows up?"

# Row 4: reworded to be specifically about making calls, assigned to Hugo
$ws.Range("A4").Value = "Debug todas as combinações de fazer chamada "
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "Hugo"
$ws.Range("D4").Value = "TBD"
$ws.Range("F4").Value = "debug"

# Row 5: still-open task, unchanged
$ws.Range("A5").Value = "Usar lista de contactos em incoming call"
$ws.Range("B5").Value = 3
$ws.Range("D5").Value = "TBD"
$ws.Range("F5").Value = "new functionality"

# Row 6: earcons work is finished -> assign + mark Done
$ws.Range("A6").Value = "Aplicar earcons"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "Hugo"
$ws.Range("D6").Value = "Done"
$ws.Range("F6").Value = "new functionality"
$ws.Range("G6").Value = "Earcons em acções de select e voltar"

# Row 7: "TTS spell function" finished and renamed -> "TTS Spell"
$ws.Range("A7").Value = "TTS Spell"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "Hugo"
$ws.Range("D7").Value = "Done"
$ws.Range("F7").Value = "new functionality"
$ws.Range("G7").Clear()

# Row 8: "Increase InCall TTS volume" replaced by handling private numbers,
# also finished
$ws.Range("A8").Value = "Lidar com numero privado (incoming number = null)"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "Hugo"
$ws.Range("D8").Value = "Done"
$ws.Range("F8").Value = "new functionality"

# Rows 9-12: already-Done tasks, now all owned by Hugo
$ws.Range("A9").Value = "Menu Manager Singleton"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = "Hugo"
$ws.Range("D9").Value = "Done"
$ws.Range("F9").Value = "Modification"
$ws.Range("G9").Value = "Reuse code"

$ws.Range("A10").Value = "No segundo ciclo de leitura não permite escolha da ultima opção"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "Hugo"
$ws.Range("D10").Value = "Done"
$ws.Range("F10").Value = "bug"

$ws.Range("A11").Value = "Allow option selection while reading title"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "Hugo"
$ws.Range("D11").Value = "Done"
$ws.Range("F11").Value = "modification"

$ws.Range("A12").Value = "Full screen"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "Hugo"
$ws.Range("D12").Value = "Done"
$ws.Range("F12").Value = "modification"

# Row 13: new task covering incoming/in-call debugging
$ws.Range("A13").Value = "Debug todas as combinações de receber chamada / em chamada"
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = "Hugo"
$ws.Range("D13").Value = "Done"

$ws.Range("A5").Select()
